# LOB1049 syllabus update (build site at 2021-01-29)
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("LOB1049")

# --- Row 4 ("Name:") previously had no B/C cells; add the English course name.
# New cells created by .Value alone pick up the rows style instead of the
# column style, so first clone the correct column formats from row 5, then
# write the text.
$ws.Range("B5:C5").Copy() | Out-Null
$ws.Range("B4:C4").PasteSpecial(-4122) | Out-Null   # xlPasteFormats
$ws.Application.CutCopyMode = $false
$ws.Range("B4").Value = 'Multivariate Statistics'
$ws.Range("C4").Value = 'Multivariate Statistics'

# --- Row 8 ("Ativação:") date bump: 01/01/2016 -> 01/01/2021.
# Assigning "01/01/2021" straight to Value would be auto-parsed into a real
# date serial by Excel. Build the literal text on a scratch cell via a
# formula (so it is never interpreted as a date), then paste-special just the
# value onto the target cells so their existing text styling is untouched.
$scratch = $ws.Range("Z1")
$scratch.Formula = "=""01/01/2021"""
$scratch.Copy() | Out-Null
$ws.Range("B8").PasteSpecial(-4163) | Out-Null   # xlPasteValues
$ws.Range("C8").PasteSpecial(-4163) | Out-Null   # xlPasteValues
$scratch.Clear() | Out-Null
$ws.Application.CutCopyMode = $false

# --- Row 14 ("Programa resumido:").
$ws.Range("B14").Value = 'Probabilidade: Vetor de variáveis aleatórias, Distribuição conjunta/marginal, Esperança e variância condicional/marginal. Estatística: Regressão Logística simples, Teste Qui-Quadrado, Testes de normalidade, Testes não-paramétricos. Técnicas Multivariadas: Gráficos multivariados, Regressão Linear Múltipla, Regressão Logística Múltipla, Análise de Variância Múltipla; Análise de agrupamento; Análise de componentes principais; Análise fatorial; Análise discriminante e Análise de correspondência'
$ws.Range("C14").Value = 'Probabilidade: Vetor de variáveis aleatórias, Distribuição conjunta/marginal, Esperança e variância condicional/marginal. Estatística: Regressão Logística simples, Teste Qui-Quadrado, Testes de normalidade, Testes não-paramétricos. Técnicas Multivariadas: Gráficos multivariados, Regressão Linear Múltipla, Regressão Logística Múltipla, Análise de Variância Múltipla; Análise de agrupamento; Análise de componentes principais; Análise fatorial; Análise discriminante e Análise de correspondência'

# --- Row 16 ("Programa:").
$ws.Range("B16").Value = 'Probabilidade: Vetor de variáveis aleatórias, Distribuição conjunta/marginal, Esperança e Variância condicional/marginal.Estatística: Regressão Logística simples (coeficiente de associação, sensitividade e especificidade, risco relativo, razão de chances), Teste Qui-Quadrado (testes de aderência, homogeneidade e independência), Teste de normalidade (Shapiro-Wilk, Teste de Kolmogorov-Smirnov), Testes não-paramétricos para amostras pareadas e independentes.Técnicas Multivariadas: Gráficos multivariados, Regressão Linear Múltipla, Regressão Logística Múltipla, Análise de Variância Múltipla; Análise de agrupamento; Análise de componentes principais; Análise fatorial; Análise discriminante e Análise de correspondência'
$ws.Range("C16").Value = 'Probabilidade: Vetor de variáveis aleatórias, Distribuição conjunta/marginal, Esperança e Variância condicional/marginal.Estatística: Regressão Logística simples (coeficiente de associação, sensitividade e especificidade, risco relativo, razão de chances), Teste Qui-Quadrado (testes de aderência, homogeneidade e independência), Teste de normalidade (Shapiro-Wilk, Teste de Kolmogorov-Smirnov), Testes não-paramétricos para amostras pareadas e independentes.Técnicas Multivariadas: Gráficos multivariados, Regressão Linear Múltipla, Regressão Logística Múltipla, Análise de Variância Múltipla; Análise de agrupamento; Análise de componentes principais; Análise fatorial; Análise discriminante e Análise de correspondência'

# --- Row 19 ("Método:").
$ws.Range("B19").Value = 'NF=A avaliação será composta por provas, listas, projetos, seminários e outras formas que farão a composição das notas, sendo estipulada a média final a somatória destas notas (N), com no mínimo duas avaliações, sendo: N = (N1+...+Nn)/n'
$ws.Range("C19").Value = 'NF=A avaliação será composta por provas, listas, projetos, seminários e outras formas que farão a composição das notas, sendo estipulada a média final a somatória destas notas (N), com no mínimo duas avaliações, sendo: N = (N1+...+Nn)/n'

# --- Row 20 ("Critério:").
$ws.Range("B20").Value = 'NF≥ 5,0.'
$ws.Range("C20").Value = 'NF≥ 5,0.'

# --- Row 21 ("Norma de recuperação:").
$ws.Range("B21").Value = '(NF+RC)/2 ≥ 5,0, onde RC é uma prova de recuperação a ser aplicada.'
$ws.Range("C21").Value = '(NF+RC)/2 ≥ 5,0, onde RC é uma prova de recuperação a ser aplicada.'

# --- Row 22 ("Bibliografia:").
$ws.Range("B22").Value = 'G.C. Runger, D. Montgomery. Estatística aplicada e probabilidade para engenheiros. São Paulo: Ed. LTC, 2009. D. C. Montgomery, E. A. Peck, G. G. Vining, Introduction to Linear Regression Analysis, 4th ed., Hoboken: John Wiley, 2006.W. J. Conover, Practical Nonparametric Statistics, 3rd ed., New York: John Wiley d Sons, 1999.R. A. Johnson, D. W. Wichern, Applied Multivariate Statistical Analysis, 6th ed., New Jersey: Prentice Hall, 2007.'
$ws.Range("C22").Value = 'G.C. Runger, D. Montgomery. Estatística aplicada e probabilidade para engenheiros. São Paulo: Ed. LTC, 2009. D. C. Montgomery, E. A. Peck, G. G. Vining, Introduction to Linear Regression Analysis, 4th ed., Hoboken: John Wiley, 2006.W. J. Conover, Practical Nonparametric Statistics, 3rd ed., New York: John Wiley d Sons, 1999.R. A. Johnson, D. W. Wichern, Applied Multivariate Statistical Analysis, 6th ed., New Jersey: Prentice Hall, 2007.'

